# Update the SMB "user name" help paragraph: change its indent from a
# first-line indent to a left (hanging-block) indent, and replace the
# short description text with the longer "local account" explanation.

$d = $word.ActiveDocument

$oldText = "Укажите имя пользователя, установленное на ПК/NAS."
$newText = "Имя учетной записи пользователя для подключения к хосту. Учетная запись Microsoft не может использоваться с SMBSync2. Создайте локальную учетную запись и используйте ее."

$found = $false
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text.TrimEnd("`r", "`a") -eq $oldText) {
        $p.Range.ParagraphFormat.FirstLineIndent = 0
        $p.Range.ParagraphFormat.LeftIndent = 21
        $r.Text = $newText
        $found = $true
        break
    }
}

if (-not $found) {
    throw "Target paragraph not found"
}
